$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value while preserving it as literal text
# (avoids Excel auto-converting numeric-looking strings to numbers,
# which would drop significant trailing/leading zeros).
function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = "28.888.14"
$ws.Range("E2").Value = "  -2.30%  "
$ws.Range("D3").Value = "1.902.25"
$ws.Range("E3").Value = "  -4.50%  "
$ws.Range("E4").Value = "  +0.25%  "
Set-TextValue $ws.Range("D5") "324.12"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  -2.03%  "
Set-TextValue $ws.Range("D8") "0.3808"
$ws.Range("E8").Value = "  -3.69%  "
Set-TextValue $ws.Range("D9") "45.56"
$ws.Range("E9").Value = "  -2.14%  "
Set-TextValue $ws.Range("D10") "0.07720"
$ws.Range("E10").Value = "  -3.01%  "
Set-TextValue $ws.Range("D11") "0.9816"
$ws.Range("E11").Value = "  -2.05%  "
Set-TextValue $ws.Range("D12") "22.04"
$ws.Range("E12").Value = "  -3.83%  "
$ws.Range("D13").Value = "1.895.40"
$ws.Range("E13").Value = "  -5.28%  "
Set-TextValue $ws.Range("D14") "6.967"
$ws.Range("E14").Value = "  -4.15%  "
$ws.Range("E15").Value = "  -3.53%  "
Set-TextValue $ws.Range("D16") "0.07048"
$ws.Range("E16").Value = "  -1.27%  "
$ws.Range("E17").Value = "  +0.16%  "
Set-TextValue $ws.Range("D18") "84.09"
$ws.Range("E18").Value = "  -5.46%  "
Set-TextValue $ws.Range("D19") "0.000009536"
$ws.Range("E19").Value = "  -4.53%  "
$ws.Range("E20").Value = "  -3.99%  "
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "28.854.90"
$ws.Range("E22").Value = "  -2.65%  "
Set-TextValue $ws.Range("D23") "5.333"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("E24").Value = "  -3.18%  "
$ws.Range("D25").Value = "2.143.72"
$ws.Range("E25").Value = "  -4.24%  "
Set-TextValue $ws.Range("D26") "2.090"
$ws.Range("E26").Value = "  -0.76%  "
Set-TextValue $ws.Range("D27") "156.88"
$ws.Range("E27").Value = "  -0.63%  "
Set-TextValue $ws.Range("D28") "19.13"
$ws.Range("E28").Value = "  -2.95%  "
Set-TextValue $ws.Range("D29") "5.588"
$ws.Range("E29").Value = "  -7.17%  "
Set-TextValue $ws.Range("D30") "117.76"
$ws.Range("E30").Value = "  -1.96%  "
Set-TextValue $ws.Range("D31") "1.841"
$ws.Range("E31").Value = "  -5.76%  "
Set-TextValue $ws.Range("D32") "0.09265"
$ws.Range("E32").Value = "  -2.07%  "
Set-TextValue $ws.Range("D33") "0.8610"
$ws.Range("E33").Value = "  -5.94%  "
Set-TextValue $ws.Range("D34") "5.095"
$ws.Range("E34").Value = "  -3.34%  "
Set-TextValue $ws.Range("D35") "1.253"
$ws.Range("E35").Value = "  -7.26%  "
$ws.Range("E36").Value = "  -5.22%  "
Set-TextValue $ws.Range("D37") "0.05699"
$ws.Range("E37").Value = "  -2.75%  "
Set-TextValue $ws.Range("D38") "1.147"
$ws.Range("E38").Value = "  -2.51%  "
$ws.Range("E39").Value = "  +0.31%  "
Set-TextValue $ws.Range("D40") "0.02033"
$ws.Range("E40").Value = "  -4.42%  "
Set-TextValue $ws.Range("D41") "7.467"
$ws.Range("E41").Value = "  -5.74%  "
Set-TextValue $ws.Range("D42") "0.5512"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("E43").Value = "  -4.25%  "
Set-TextValue $ws.Range("D44") "9.301"
$ws.Range("E44").Value = "  -5.82%  "
Set-TextValue $ws.Range("D45") "2.730"
$ws.Range("E45").Value = "  -0.65%  "
Set-TextValue $ws.Range("D46") "0.5193"
$ws.Range("E46").Value = "  -3.76%  "
$ws.Range("E47").Value = "  -6.89%  "
Set-TextValue $ws.Range("D48") "2.088"
$ws.Range("E48").Value = "  -4.56%  "
Set-TextValue $ws.Range("D49") "0.06824"
$ws.Range("E49").Value = "  -1.87%  "
Set-TextValue $ws.Range("D50") "111.33"
$ws.Range("E50").Value = "  -2.70%  "
$ws.Range("E51").Value = "  -5.55%  "
